# This script applies a re-sort of the weekly price rows (rows 2-22) by
# permuting the "variable" columns (D, M, N, O, P, Q, R, S, T) across rows,
# while leaving the "constant" descriptive columns (A, B, C, E, F, G, H, I,
# J, K, L) untouched (they are identical for every row anyway).
#
# The mapping below says: new row <key> gets the old values that currently
# live in row <value>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 5
    3  = 14
    4  = 11
    5  = 12
    6  = 13
    7  = 18
    8  = 7
    9  = 17
    10 = 9
    11 = 20
    12 = 15
    13 = 2
    14 = 16
    15 = 6
    16 = 22
    17 = 3
    18 = 4
    19 = 10
    20 = 19
    21 = 21
    22 = 8
}

# Columns that get shuffled between rows according to the mapping above.
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current ("before") values for every relevant column/row so
# that overwriting some rows doesn't clobber data still needed for others.
$snapshot = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 22; $row++) {
        $snapshot["$col$row"] = $ws.Range("$col$row").Value2
    }
}

# Write the permuted values back.
foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $snapshot["$col$srcRow"]
    }
}
